# Update the "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
# F3: 7390 -> 7391
# F5: 90   -> 91
# F15: 475 -> 476

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 7391
    $ws.Range("F5").Value = 91
    $ws.Range("F15").Value = 476
}
